# Append 4 new report rows (32-35) to Sheet1, mirroring the existing
# "Phone Number / Status Web / Status Api / Message / Timestamp" layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("855764041418", "success", "failed", "Message not found or Archived for another partner", "2024-07-26 09:09:51"),
    @("855764045818", "success", "failed", "Message not found or Archived for another partner", "2024-07-26 09:11:31"),
    @("855764041418", "success", "failed", "Message not found or Archived for another partner", "2024-07-26 09:09:51"),
    @("855764045818", "success", "failed", "Message not found or Archived for another partner", "2024-07-26 09:11:31")
)

$startRow = 32
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]

    # Column A holds a phone-number-like string; prefix with an apostrophe so
    # Excel stores it as text (matches the other rows' Phone Number column)
    # instead of silently coercing it to a number.
    $ws.Range("A$r").Value = "'" + $values[0]
    $ws.Range("B$r").Value = $values[1]
    $ws.Range("C$r").Value = $values[2]
    $ws.Range("D$r").Value = $values[3]
    $ws.Range("E$r").Value = $values[4]
}
